# Tripadvisor New Orleans shard 94 update:
#  1. hotel_info: insert a new "State" column between Hotel_Name and City,
#     populated with "Louisiana" for the existing data row.
#  2. Reorder worksheet tabs so review_info precedes hotel_info.

$wb = $excel.ActiveWorkbook

# --- 1. Add State column to hotel_info -------------------------------------
$hotelInfo = $wb.Worksheets.Item("hotel_info")

$hotelInfo.Columns.Item(3).Insert()
$hotelInfo.Cells.Item(1, 3).Value = "State"
$hotelInfo.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Move review_info so it becomes the first sheet tab ------------------
$reviewInfo = $wb.Worksheets.Item("review_info")
$reviewInfo.Move($wb.Worksheets.Item(1))
